# Change highlight color from yellow to bright green on the four
# "friend request" related user-story paragraphs (both the run text and
# the paragraph-mark formatting), matching the canonical OOXML diff.

$d = $word.ActiveDocument

# wdBrightGreen = 4 (maps to OOXML w:highlight w:val="green")
$wdBrightGreen = 4

$targets = @(
    "As a social media user, I want to send a friend request to another user.",
    "As a social media user, I want to accept or deny a friend request from another user.",
    "As a social media user, I want to see a list of pending friend requests.",
    "As a social media user, I want to remove a friend from my list of friends."
)

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    foreach ($t in $targets) {
        if ($text -like "*$t*") {
            # Setting Font.HighlightColorIndex (rather than
            # Range.HighlightColorIndex) updates both the run's
            # highlight and the paragraph mark's highlight, matching
            # the target XML where both <w:pPr><w:rPr> and <w:r><w:rPr>
            # carry the new color.
            $p.Range.Font.HighlightColorIndex = $wdBrightGreen
        }
    }
}

Write-Output "done"
